$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 55429
$ws.Range("B28").Value = 437
$ws.Range("B71").Value = 7807
$ws.Range("B98").Value = 2499
$ws.Range("B118").Value = 10
$ws.Range("B233").Value = 66361
$ws.Range("B250").Value = 420164
$ws.Range("B251").Value = 261361
$ws.Range("B269").Value = 79532
$ws.Range("B270").Value = 160719
$ws.Range("B274").Value = 578
$ws.Range("B288").Value = 3628
$ws.Range("B297").Value = 24332
$ws.Range("B309").Value = 3359
$ws.Range("B389").Value = 8
$ws.Range("B468").Value = 4228
$ws.Range("B497").Value = 142
$ws.Range("B521").Value = 20262
$ws.Range("B566").Value = 1750
$ws.Range("B588").Value = 2
$ws.Range("B640").Value = 2691
$ws.Range("B664").Value = 13304
$ws.Range("B688").Value = 33
$ws.Range("B698").Value = 3221
$ws.Range("B710").Value = 1321
$ws.Range("B735").Value = 429
$ws.Range("B738").Value = 15663
$ws.Range("B751").Value = 5484
$ws.Range("B753").Value = 208899
$ws.Range("B773").Value = 35
$ws.Range("B781").Value = 518
$ws.Range("B786").Value = 413
$ws.Range("B797").Value = 484
$ws.Range("B804").Value = 2025
$ws.Range("B806").Value = 178
$ws.Range("B834").Value = 1876
$ws.Range("B839").Value = 153
$ws.Range("B862").Value = 4684
$ws.Range("B871").Value = 4041
$ws.Range("B881").Value = 664
$ws.Range("B895").Value = 313
$ws.Range("B904").Value = 777
$ws.Range("B922").Value = 659
$ws.Range("B938").Value = 14555
$ws.Range("B946").Value = 9738
$ws.Range("B957").Value = 5556
$ws.Range("B958").Value = 2456
$ws.Range("B979").Value = 62199
$ws.Range("B1038").Value = 31927
$ws.Range("B1043").Value = 50290
$ws.Range("B1053").Value = 579
$ws.Range("B1064").Value = 8
$ws.Range("B1065").Value = 102
$ws.Range("B1072").Value = 33
$ws.Range("B1085").Value = 44
